$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The project was renamed from "TemplateGO" to "XlsxTemplate", so the title
# banner cell in the test workbook ("TemplateGO 测试") needs to read
# "XlsxTemplate 测试" instead.
$titleCell = $ws.Cells.Find("TemplateGO")
if ($titleCell) {
    $titleCell.Value = "XlsxTemplate 测试"
} else {
    $ws.Range("A1").Value = "XlsxTemplate 测试"
}

# Row 12 is a trailing, completely empty (style-only) row below the data
# table - clean it up so the sheet's used range goes back to A1:E11.
$lastRow = $ws.Range("A12:E12")
if ($ws.Application.WorksheetFunction.CountA($lastRow) -eq 0) {
    $lastRow.EntireRow.Delete()
}
